$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.76
$ws.Range("L2").Value = 1.38
$ws.Range("O2").Value = 1.3
$ws.Range("Q2").Value = 1.87
$ws.Range("S2").Value = 3.2
$ws.Range("T2").Value = 1.83
$ws.Range("U2").Value = 2.1
$ws.Range("AI2").Value = 75
$ws.Range("L3").Value = 1.32
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 2.04
$ws.Range("Q3").Value = 1.86
$ws.Range("R3").Value = 1.41
$ws.Range("S3").Value = 3.15
$ws.Range("T3").Value = 1.72
$ws.Range("U3").Value = 2.22
$ws.Range("X3").Value = 16.5
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 27
$ws.Range("AA3").Value = 65
$ws.Range("AB3").Value = 11
$ws.Range("AC3").Value = 8.6
$ws.Range("AE3").Value = 46
$ws.Range("AF3").Value = 16.5
$ws.Range("AG3").Value = 12
$ws.Range("AI3").Value = 48
$ws.Range("AK3").Value = 24
$ws.Range("AM3").Value = 100
$ws.Range("AN3").Value = 19.5
$ws.Range("AO3").Value = 38
$ws.Range("Q5").Value = 1.88
$ws.Range("F6").Value = 2.06
$ws.Range("N6").Value = 5
$ws.Range("S6").Value = 2.54
$ws.Range("T6").Value = 1.62
$ws.Range("U6").Value = 2.38
$ws.Range("F7").Value = 2.64
$ws.Range("G7").Value = 16
$ws.Range("H7").Value = 1.38
$ws.Range("J7").Value = 4.7
$ws.Range("F8").Value = 2.04
$ws.Range("G8").Value = 2.08
$ws.Range("F9").Value = 1.9
$ws.Range("G9").Value = 2.04
$ws.Range("H9").Value = 4.3
$ws.Range("I9").Value = 5.4
$ws.Range("P9").Value = 1.76
$ws.Range("Q9").Value = 1.93
$ws.Range("G10").Value = 3.25
$ws.Range("I10").Value = 2.98
$ws.Range("J10").Value = 3.25
$ws.Range("Q10").Value = 1.6
$ws.Range("F12").Value = 1.99
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 4.6
$ws.Range("G13").Value = 1.17
$ws.Range("H13").Value = 30
$ws.Range("J13").Value = 9.4
$ws.Range("F14").Value = 1.94
$ws.Range("G14").Value = 2.5
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 4.9
$ws.Range("J14").Value = 3.1
$ws.Range("K14").Value = 5.7
$ws.Range("P14").Value = 1.61
$ws.Range("Q14").Value = 1.97
$ws.Range("F15").Value = 1.87
$ws.Range("G15").Value = 2.08
$ws.Range("I15").Value = 4.9
$ws.Range("J15").Value = 3.85
$ws.Range("H16").Value = 1.91
$ws.Range("J16").Value = 3.5
$ws.Range("P16").Value = 1.86
$ws.Range("S16").Value = 3.55
$ws.Range("G17").Value = 2.14
$ws.Range("K17").Value = 3.55
$ws.Range("N17").Value = 3.25
$ws.Range("T17").Value = 1.86
$ws.Range("AB17").Value = 8.6
$ws.Range("F18").Value = 1.9
$ws.Range("G18").Value = 1.99
$ws.Range("H18").Value = 3.95
$ws.Range("I18").Value = 4.4
$ws.Range("J18").Value = 4
$ws.Range("N18").Value = 3.9
$ws.Range("O18").Value = 1.29
$ws.Range("Q18").Value = 1.83
$ws.Range("U18").Value = 2.06
$ws.Range("X18").Value = 21
$ws.Range("AH18").Value = 24
$ws.Range("AM18").Value = 130
$ws.Range("F19").Value = 2.54
$ws.Range("H19").Value = 2.88
$ws.Range("Q19").Value = 1.92
$ws.Range("AA19").Value = 50
$ws.Range("G20").Value = 1.56
$ws.Range("S20").Value = 2.92
$ws.Range("T20").Value = 1.9
$ws.Range("AI21").Value = 48
$ws.Range("AM21").Value = 95
$ws.Range("AN21").Value = 25
$ws.Range("AO21").Value = 30
$ws.Range("F22").Value = 1.35
$ws.Range("H22").Value = 9.4
$ws.Range("K22").Value = 5.9
$ws.Range("F23").Value = 2.4
$ws.Range("G23").Value = 2.62
$ws.Range("H23").Value = 2.7
$ws.Range("P23").Value = 2.38
$ws.Range("H24").Value = 9.6
$ws.Range("J24").Value = 5.2
$ws.Range("N24").Value = 4.6
$ws.Range("P24").Value = 2.24
$ws.Range("F25").Value = 1.71
$ws.Range("G25").Value = 1.91
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 4.2
$ws.Range("H26").Value = 4.5
$ws.Range("J26").Value = 3.9
$ws.Range("F27").Value = 2.4
$ws.Range("G27").Value = 2.72
$ws.Range("H27").Value = 2.6
$ws.Range("I27").Value = 2.98
$ws.Range("F28").Value = 2.3
$ws.Range("G28").Value = 2.6
$ws.Range("H28").Value = 3.05
$ws.Range("I28").Value = 4
$ws.Range("J28").Value = 3.2
$ws.Range("K28").Value = 4.3
$ws.Range("P28").Value = 1.78
$ws.Range("Q28").Value = 2
$ws.Range("F29").Value = 2.84
$ws.Range("G29").Value = 3.25
$ws.Range("I29").Value = 2.64
$ws.Range("P29").Value = 2.04
$ws.Range("Q29").Value = 1.76
$ws.Range("P30").Value = 2.2
$ws.Range("Q30").Value = 1.66
$ws.Range("I31").Value = 11.5
$ws.Range("Q31").Value = 1.65
$ws.Range("G32").Value = 3.2
$ws.Range("H32").Value = 2.66
$ws.Range("I32").Value = 3.35
$ws.Range("J32").Value = 3.1
$ws.Range("K32").Value = 4.3
$ws.Range("P32").Value = 1.92
$ws.Range("G33").Value = 1.69
$ws.Range("I33").Value = 7.6
$ws.Range("J33").Value = 4.3
$ws.Range("Q33").Value = 1.51
$ws.Range("O34").Value = 1.53
$ws.Range("S34").Value = 5.3
$ws.Range("F35").Value = 1.78
$ws.Range("G35").Value = 1.89
$ws.Range("H35").Value = 4.6
$ws.Range("I35").Value = 5.2
$ws.Range("K35").Value = 4.3
$ws.Range("P35").Value = 2.08
$ws.Range("Q35").Value = 1.67
$ws.Range("F36").Value = 1.54
$ws.Range("G36").Value = 1.65
$ws.Range("P36").Value = 1.89
$ws.Range("Q36").Value = 1.93
$ws.Range("F37").Value = 2.86
$ws.Range("G37").Value = 38
$ws.Range("H37").Value = 1.35
$ws.Range("I37").Value = 1.53
$ws.Range("J37").Value = 4.3
$ws.Range("P37").Value = 1.92
$ws.Range("Q37").Value = 1.65
$ws.Range("G38").Value = 2.32
$ws.Range("H38").Value = 3.6
$ws.Range("P38").Value = 1.84
$ws.Range("Q38").Value = 1.87
$ws.Range("R39").Value = 1.6
$ws.Range("S39").Value = 2.56
$ws.Range("U39").Value = 2.54
$ws.Range("F41").Value = 1.64
$ws.Range("H41").Value = 1.99
$ws.Range("I41").Value = 9.199999999999999
$ws.Range("J41").Value = 3.6
$ws.Range("P41").Value = 1.9
$ws.Range("Q41").Value = 1.71
$ws.Range("J42").Value = 3.4
$ws.Range("F43").Value = 7.4
$ws.Range("G43").Value = 13.5
$ws.Range("J43").Value = 5.1
$ws.Range("K43").Value = 6.2
$ws.Range("AE44").Value = 500
$ws.Range("Q46").Value = 1.7
$ws.Range("Z46").Value = 9.199999999999999
$ws.Range("AL46").Value = 85
$ws.Range("F47").Value = 1.76
$ws.Range("G47").Value = 1.95
$ws.Range("H47").Value = 4.4
$ws.Range("I47").Value = 5.3
$ws.Range("J47").Value = 3.65
$ws.Range("K47").Value = 4.3
$ws.Range("P47").Value = 1.99
$ws.Range("Q47").Value = 1.8
$ws.Range("S48").Value = 3.05
$ws.Range("U48").Value = 2.36
$ws.Range("F49").Value = 2.38
$ws.Range("H49").Value = 3.25
$ws.Range("I49").Value = 3.5
$ws.Range("J49").Value = 3.2
$ws.Range("P49").Value = 2.08
$ws.Range("Q49").Value = 1.81
$ws.Range("P50").Value = 2.12
$ws.Range("F51").Value = 1.23
$ws.Range("H51").Value = 12
$ws.Range("I51").Value = 30
$ws.Range("K51").Value = 7.8
$ws.Range("Q51").Value = 1.52
$ws.Range("P52").Value = 1.93
$ws.Range("Q52").Value = 1.94
$ws.Range("F54").Value = 6
$ws.Range("U54").Value = 2.64
$ws.Range("F55").Value = 2.12
$ws.Range("I55").Value = 4.4
$ws.Range("J55").Value = 3.2
$ws.Range("S56").Value = 2.82
$ws.Range("G57").Value = 3.95
$ws.Range("H57").Value = 2.26
